# Auto-generated Excel COM-interop script
# Applies the 2025-12-23 daily crime count increments to the violent-crime-full-year workbook.
# For every affected worksheet, update the specific cells (column L = 2025 running totals,
# and in two cases column J = 2023 totals for the "NO NEIGHBORHOOD DATA" bucket correction)
# to their new values as captured in the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 6457
$ws.Range("J3").Value = 8083
$ws.Range("L3").Value = 6961
$ws.Range("L4").Value = 1736
$ws.Range("L5").Value = 410
$ws.Range("L6").Value = 5717
$ws.Range("J7").Value = 29363
$ws.Range("L7").Value = 21281

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 189
$ws.Range("L4").Value = 77
$ws.Range("L7").Value = 674
$ws.Range("L8").Value = 1407
$ws.Range("L11").Value = 349
$ws.Range("L15").Value = 180
$ws.Range("L18").Value = 146
$ws.Range("L19").Value = 585
$ws.Range("L23").Value = 223
$ws.Range("L27").Value = 183
$ws.Range("L29").Value = 1191
$ws.Range("L31").Value = 213
$ws.Range("L33").Value = 958
$ws.Range("L37").Value = 817
$ws.Range("L41").Value = 91
$ws.Range("L42").Value = 671
$ws.Range("L43").Value = 158
$ws.Range("L49").Value = 114
$ws.Range("L52").Value = 451
$ws.Range("L53").Value = 241
$ws.Range("L54").Value = 458
$ws.Range("L55").Value = 226
$ws.Range("L56").Value = 22
$ws.Range("L58").Value = 12
$ws.Range("J63").Value = 239
$ws.Range("L63").Value = 67
$ws.Range("L65").Value = 419
$ws.Range("L67").Value = 737
$ws.Range("L69").Value = 47
$ws.Range("L73").Value = 168
$ws.Range("L76").Value = 335
$ws.Range("L78").Value = 283
$ws.Range("L79").Value = 589
$ws.Range("L80").Value = 69
$ws.Range("L83").Value = 468
$ws.Range("L84").Value = 203
$ws.Range("L85").Value = 1054
$ws.Range("L89").Value = 287
$ws.Range("L91").Value = 287
$ws.Range("L93").Value = 106
$ws.Range("L96").Value = 234
$ws.Range("J101").Value = 29363
$ws.Range("L101").Value = 21281

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L4").Value = 25
$ws.Range("L7").Value = 234

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L6").Value = 160
$ws.Range("L7").Value = 674

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 98
$ws.Range("L7").Value = 349

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L2").Value = 74
$ws.Range("L7").Value = 287

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L6").Value = 219
$ws.Range("L7").Value = 1054

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L5").Value = 13
$ws.Range("L7").Value = 451

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 47

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 79
$ws.Range("L7").Value = 241

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 433
$ws.Range("L3").Value = 495
$ws.Range("L6").Value = 338
$ws.Range("L7").Value = 1407

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 150
$ws.Range("L7").Value = 468

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L6").Value = 269
$ws.Range("L7").Value = 958

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L4").Value = 49
$ws.Range("L6").Value = 209
$ws.Range("L7").Value = 817

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 155
$ws.Range("L7").Value = 419

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L6").Value = 56
$ws.Range("L7").Value = 213

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 287
$ws.Range("L6").Value = 171
$ws.Range("L7").Value = 737

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L2").Value = 70
$ws.Range("L7").Value = 203

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L6").Value = 46
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 218
$ws.Range("L7").Value = 458

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 458
$ws.Range("L7").Value = 1191

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 209
$ws.Range("L3").Value = 178
$ws.Range("L6").Value = 161
$ws.Range("L7").Value = 585

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L3").Value = 66
$ws.Range("L7").Value = 335

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L2").Value = 28
$ws.Range("L7").Value = 91

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 183
$ws.Range("L3").Value = 230
$ws.Range("L7").Value = 671

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L4").Value = 33
$ws.Range("L6").Value = 82
$ws.Range("L7").Value = 283

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L3").Value = 75
$ws.Range("L7").Value = 226

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L4").Value = 19
$ws.Range("L6").Value = 56
$ws.Range("L7").Value = 223

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L6").Value = 38
$ws.Range("L7").Value = 287

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 191
$ws.Range("L7").Value = 589

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L4").Value = 16
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 146

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L2").Value = 37
$ws.Range("L7").Value = 106

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L4").Value = 18
$ws.Range("L7").Value = 180

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L6").Value = 39
$ws.Range("L7").Value = 168

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L3").Value = 62
$ws.Range("L7").Value = 189

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L3").Value = 51
$ws.Range("L7").Value = 183

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L6").Value = 48
$ws.Range("L7").Value = 158

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 22

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item("Millenium Park")
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 12
